$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cells below can be edited.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A10):
# 2021-07-09 -> 2021-07-13
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) / Percent Change (E) figures for each holding (rows 2-7)
$ws.Range("D2").Value = 0.4965621713973384
$ws.Range("E2").Value = -0.004209720627631031

$ws.Range("D3").Value = 0.3250290934955451
$ws.Range("E3").Value = -0.007099583491101846

$ws.Range("D4").Value = 0.09182021818462621
$ws.Range("E4").Value = -0.01461892247043362

$ws.Range("D5").Value = 0.05837823293739635
$ws.Range("E5").Value = -0.002955553029441838

$ws.Range("D6").Value = 0.028210283985094
$ws.Range("E6").Value = -0.0205726994717822

$ws.Range("E7").Value = -0.006493173510651373

# Restore sheet protection to match the original, published state.
$ws.Protect()
